# Rename the "DimAssinatura" table's 5th column header (cell E1) from
# " Mensal " (with leading/trailing spaces) to "Mensal" (trimmed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Mensal"

# Update the active selection to E1, matching the saved view state.
$ws.Range("E1").Select()
